$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "242.48"
Set-TextValue "D3" "21.54"
Set-TextValue "D4" "5.245"
Set-TextValue "D5" "0.05607"
Set-TextValue "D6" "3.373"
Set-TextValue "D7" "6.379"
Set-TextValue "D8" "0.8078"
Set-TextValue "D9" "0.9149"
Set-TextValue "D10" "0.1427"
Set-TextValue "D11" "0.07313"
Set-TextValue "D12" "0.03214"
Set-TextValue "D13" "0.03012"
Set-TextValue "D15" "3.614"
Set-TextValue "D16" "0.001644"
Set-TextValue "D17" "0.04705"
Set-TextValue "D18" "0.0005816"
Set-TextValue "D19" "0.006353"
Set-TextValue "D20" "0.004977"
Set-TextValue "D21" "0.001043"
Set-TextValue "D22" "0.0001502"
Set-TextValue "D23" "0.0003103"
Set-TextValue "D24" "3.766"
Set-TextValue "D26" "0.3271"
Set-TextValue "D40" "0.03913"
Set-TextValue "D41" "0.006966"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003403"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1033"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.007498"
Set-TextValue "D45" "0.00005935"
Set-TextValue "D47" "0.0005506"
Set-TextValue "D48" "0.6831"
Set-TextValue "D49" "0.06309"
Set-TextValue "D50" "0.00002102"
